$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in the "+" action-button cells in column O for rows 7 and 8,
# matching the formatting already used for the rest of column O (rows 2-6)
# and the adjacent button cells on the same rows.
$ws.Range("O6").Copy()
$ws.Range("O7").PasteSpecial(-4104)

$ws.Range("M8").Copy()
$ws.Range("O8").PasteSpecial(-4104)

$excel.CutCopyMode = 0

# Move the active selection to O8, matching the saved view state.
$ws.Range("O8").Select()
